$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove row 8 entirely, then insert a blank row back in its place so that
# rows below keep their original row numbers (e.g. row 12 stays row 12)
$ws.Rows("8:8").Delete()
$ws.Rows("8:8").Insert()
$ws.Rows("8:8").Clear()

# Update selected cell shown in the sheet view
$ws.Range("M2").Select()
